# Entire flow from invoice submission to view payment request done for po
# and support for multiple browsers added in test base.
#
# This script updates the two "memo" test-data sheets (Memo_Verification_details
# and Memo_invoice_Details) with a new sample invoice (TESTINV90008) replacing
# the previous sample invoice (TESTINV45388), and updates the selected cell on
# the verification sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Memo_Verification_details")
$ws2 = $wb.Worksheets.Item("Memo_invoice_Details")

# --- Memo_Verification_details (sheet1), row 2 ---------------------------
$ws1.Range("B2").Value = "TESTINV90008"   # Invoice_number
$ws1.Range("K2").Value = "V0"             # Revised Tax Code
$ws1.Range("O2").Value = "test"           # Assignment

# Move the active selection from K2 to N2
$ws1.Activate()
$ws1.Range("N2").Select()

# --- Memo_invoice_Details (sheet2), row 2 ---------------------------------
$ws2.Range("B2").Value = "TESTINV90008"   # Invoice_number

# Invoice_Date must stay a literal text value (it looks like a date, but the
# source sheet stores it as plain text) - force text entry so Excel doesn't
# silently convert it to a date serial number.
$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Value = "2024-03-13"     # Invoice_Date
$ws2.Range("C2").NumberFormat = "General"

# IGST / TD / TCS / Total_Invoice_Amount are stored as text (not numbers) in
# this sheet, same as the rest of the row, so force text entry the same way
# as the Invoice_Date above (otherwise Excel auto-converts "0"/"1" to real
# numbers).
$ws2.Range("I2").NumberFormat = "@"
$ws2.Range("I2").Value = "0"              # IGST
$ws2.Range("I2").NumberFormat = "General"
$ws2.Range("I2").HorizontalAlignment = -4152  # xlRight

$ws2.Range("J2").NumberFormat = "@"
$ws2.Range("J2").Value = "0"              # TD
$ws2.Range("J2").NumberFormat = "General"

$ws2.Range("K2").NumberFormat = "@"
$ws2.Range("K2").Value = "0"              # TCS
$ws2.Range("K2").NumberFormat = "General"

$ws2.Range("M2").Value = "test"           # Customer_Name
$ws2.Range("N2").Value = "test"           # Comments

$ws2.Range("O2").NumberFormat = "@"
$ws2.Range("O2").Value = "1"              # Total_Invoice_Amount
$ws2.Range("O2").NumberFormat = "General"
$ws2.Range("O2").HorizontalAlignment = -4152  # xlRight
